# Rename the product "LenderCoin" -> "LendCoin" (and "LenderCoins" -> "LendCoins")
# throughout the deck. The term appears as its own run (flagged err="1" by the
# spell-checker) inside larger title/body paragraphs on a few slides.

$p = $ppt.ActivePresentation

function Replace-TextOnSlide($slideIndex, $searchText, $replaceText) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            $found = $tr.Find($searchText, 0, $false, $false)
            while ($found -ne $null) {
                $found.Text = $replaceText
                $found = $tr.Find($searchText, 0, $false, $false)
            }
        }
    }
}

# Slide 3 title: "Our Proposal - LenderCoin" -> "Our Proposal - LendCoin"
Replace-TextOnSlide 3 "LenderCoin" "LendCoin"

# Slide 4 title: "About LenderCoin" -> "About LendCoin"
Replace-TextOnSlide 4 "LenderCoin" "LendCoin"

# Slide 5 body: "Each user of LenderCoin is given..." -> "...LendCoin..."
Replace-TextOnSlide 5 "LenderCoin" "LendCoin"

# Slide 8 body: "...awarded a payment of 100 LenderCoins." -> "...LendCoins."
# (handled with the plural form so the singular replacement above doesn't
# clobber it first)
Replace-TextOnSlide 8 "LenderCoins" "LendCoins"
